$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" (Changed) column C was refreshed from 2023-10-05 (45204)
# to 2023-10-06 (45205) for every data row (rows 2 through 237).
$ws.Range("C2:C237").Value = 45205
